$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. ETI sheet: fill in the newly-entered survey rows (B:O, rows 2-8)
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ETI")
$ws.Select()

# numeric answer columns (B..L), row by row
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1

$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 4

$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 4

$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 4

$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 3

$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 3
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 2

$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 2

# essay columns (M, N, O) - written column-by-column so the new
# shared-string entries land in the same order as the source workbook
$ws.Range("M2").Value = "asd asd asd"
$ws.Range("M3").Value = "dsa dsa dsa"
$ws.Range("M4").Value = "sdf sdf sdf"
$ws.Range("M5").Value = "fds fds fds"
$ws.Range("M6").Value = "dfg dfg dfg"
$ws.Range("M7").Value = "gfd gfd gfd"
$ws.Range("M8").Value = "fgh fgh fgh"

$ws.Range("N2").Value = "qwe qwe qwe"
$ws.Range("N3").Value = "ewq ewq ewq"
$ws.Range("N4").Value = "wer wer wer"
$ws.Range("N5").Value = "rew rew rew"
$ws.Range("N6").Value = "ert ert ert"
$ws.Range("N7").Value = "tre tre tre"
$ws.Range("N8").Value = "rty rty rty"

$ws.Range("O2").Value = "zxc zxc zxc"
$ws.Range("O3").Value = "cxz cxz cxz"
$ws.Range("O4").Value = "xcv xcv xcv"
$ws.Range("O5").Value = "vcx vcx vcx"
$ws.Range("O6").Value = "cvb cvb cvb"
$ws.Range("O7").Value = "bvc bvc bvc"
$ws.Range("O8").Value = "vbn vbn vbn"

# update the view selection on the ETI sheet
$ws.Range("F6").Select()

# ------------------------------------------------------------------
# 2. trainer sheet: move the selection from A3 to B3
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("trainer")
$ws2.Select()
$ws2.Range("B3").Select()

# ------------------------------------------------------------------
# 3. main sheet: move the selection from K7 to I24 (and keep it the
#    active/displayed sheet, matching the original workbook state)
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("main")
$ws1.Select()
$ws1.Range("I24").Select()
